$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.948.98"
$ws.Range("E2").Value = "  -1.46%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.891.71"
$ws.Range("E3").Value = "  -2.32%  "

# Row 5 - XRP
$ws.Range("D5").Value = "'0.7349"
$ws.Range("E5").Value = "  -1.38%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'242.52"
$ws.Range("E6").Value = "  -1.13%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.30%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3085"
$ws.Range("E8").Value = "  -2.67%  "

# Row 9 - Solana
$ws.Range("D9").Value = "'26.34"
$ws.Range("E9").Value = "  -4.15%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.06889"
$ws.Range("E10").Value = "  -1.33%  "

# Row 11 - now TRON (was Polygon)
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.07954"
$ws.Range("E11").Value = "  -0.56%  "

# Row 12 - now Polygon (was TRON)
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7687"
$ws.Range("E12").Value = "  -1.35%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.889.33"
$ws.Range("E13").Value = "  -2.48%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.207"

# Row 15 - Litecoin
$ws.Range("D15").Value = "'91.44"
$ws.Range("E15").Value = "  -3.08%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "29.960.99"

# Row 17 - Avalanche
$ws.Range("D17").Value = "'14.07"
$ws.Range("E17").Value = "  -2.31%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "'5.795"
$ws.Range("E18").Value = "  +1.16%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'239.41"
$ws.Range("E19").Value = "  -5.18%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "'0.000007754"
$ws.Range("E20").Value = "  -2.26%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.21%  "

# Row 22 - WrappedliquidstakedEther2.0
$ws.Range("D22").Value = "2.151.81"
$ws.Range("E22").Value = "  -1.87%  "

# Row 23 - BinanceUSD
$ws.Range("E23").Value = "  -0.25%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "'6.946"
$ws.Range("E24").Value = "  +4.13%  "

# Row 25 - Cosmos
$ws.Range("D25").Value = "'9.272"
$ws.Range("E25").Value = "  -2.16%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'165.09"
$ws.Range("E26").Value = "  -0.61%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  -1.00%  "

# Row 28 - Stellar
$ws.Range("D28").Value = "'0.1265"
$ws.Range("E28").Value = "  -4.20%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'2.020"
$ws.Range("E29").Value = "  -9.57%  "

# Row 30 - Toncoin
$ws.Range("D30").Value = "'1.350"
$ws.Range("E30").Value = "  -1.04%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "'1.531"
$ws.Range("E31").Value = "  +1.38%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'4.290"
$ws.Range("E32").Value = "  -1.56%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "'4.053"
$ws.Range("E33").Value = "  -1.10%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "'0.05090"
$ws.Range("E34").Value = "  -1.27%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  +0.21%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "'0.7338"
$ws.Range("E36").Value = "  -1.55%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value = "'2.721"
$ws.Range("E37").Value = "  -2.36%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "'0.01920"
$ws.Range("E38").Value = "  -1.11%  "

# Row 39 - MXToken
$ws.Range("D39").Value = "'2.783"
$ws.Range("E39").Value = "  -0.67%  "

# Row 40 - FraxShare
$ws.Range("D40").Value = "'6.292"
$ws.Range("E40").Value = "  -1.93%  "

# Row 41 - Aave
$ws.Range("D41").Value = "'74.08"
$ws.Range("E41").Value = "  -4.62%  "

# Row 42 - TheSandbox
$ws.Range("D42").Value = "'0.4432"
$ws.Range("E42").Value = "  -0.55%  "

# Row 43 - RenderToken
$ws.Range("E43").Value = "  -1.62%  "

# Row 44 - PaxDollar
$ws.Range("E44").Value = "  -0.24%  "

# Row 45 - TrustWalletToken
$ws.Range("D45").Value = "'0.8344"
$ws.Range("E45").Value = "  +0.30%  "

# Row 46 - Aptos
$ws.Range("D46").Value = "'7.626"
$ws.Range("E46").Value = "  +2.43%  "

# Row 47 - Quant
$ws.Range("D47").Value = "'100.83"
$ws.Range("E47").Value = "  -0.47%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "'9.781"
$ws.Range("E48").Value = "  +0.60%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "2.049.60"
$ws.Range("E49").Value = "  -1.82%  "

# Row 50 - Elrond
$ws.Range("D50").Value = "'36.35"
$ws.Range("E50").Value = "  -2.26%  "

# Row 51 - Maker
$ws.Range("D51").Value = "'931.78"
$ws.Range("E51").Value = "  -5.51%  "
